# New crime data collected - weekly CompStat update for the 123rd Precinct.
# Updates: report volume/date header text, several Week-to-Date / 28-Day /
# Year-to-Date crime-count cells (and their derived % change columns) for
# rows 15-21, 24-25, and 26-30, plus the autosized width of column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 30   Number  44" -> "...  45"  (rich-text run, in place)
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "45"

# ---------------------------------------------------------------------
# Header: "Report Covering the Week  10/30/2023  Through  11/5/2023"
#      -> "Report Covering the Week  11/6/2023  Through  11/12/2023"
# ---------------------------------------------------------------------
$ws.Range("C9").Characters(27, 10).Text = "11/6/2023"
$ws.Range("C9").Characters(47, 9).Text = "11/12/2023"

# ---------------------------------------------------------------------
# Column E autofit width grew slightly because of the new values below.
# ---------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 8.699091

# Number formats used throughout the crime table.
$fmtGeneral = "General"
$fmtText = "@"
$fmtCount = "#,##0"
$fmtPct1 = '#,##0.0;"-"#,##0.0'

function Set-TextZero($addr) {
    # Turn a numeric cell into the literal text "0" (style used across the
    # sheet for "no data this period"), keeping the sheet's normal font.
    $c = $ws.Range($addr)
    $c.NumberFormat = $fmtText
    $c.Value = "0"
    $c.NumberFormat = $fmtGeneral
}

function Set-TextStar($addr) {
    # Turn a numeric cell into the literal text "***.*" (the sheet's
    # placeholder for an undefined % change).
    $c = $ws.Range($addr)
    $c.NumberFormat = $fmtText
    $c.Value = "***.*"
    $c.NumberFormat = $fmtGeneral
}

function Set-Count($addr, $value) {
    # Turn a "0"/General text cell into a real count number.
    $c = $ws.Range($addr)
    $c.Value = $value
    $c.NumberFormat = $fmtCount
}

function Set-Pct($addr, $value, $fmt) {
    # Turn a "***.*"/General text cell into a real % change number.
    $c = $ws.Range($addr)
    $c.Value = $value
    $c.NumberFormat = $fmt
}

# ---------------------------------------------------------------------
# Row 15 - Rape: Week-to-date 2023 count drops to "no data" this week.
# ---------------------------------------------------------------------
Set-TextZero "D15"
Set-TextStar "E15"

# ---------------------------------------------------------------------
# Row 16 - Robbery: Week-to-date 2023 now has 1 complaint (was "no data").
# ---------------------------------------------------------------------
Set-Count "C16" 1
$ws.Range("I16").Value = 18
$ws.Range("K16").Value = 5.882352941176
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -18.181818181818
$ws.Range("N16").Value = -68.965517241379

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault.
# ---------------------------------------------------------------------
Set-Count "C17" 1
Set-TextZero "D17"
Set-TextStar "E17"
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 77
$ws.Range("K17").Value = 79.069767441860
$ws.Range("L17").Value = 156.666666666667
$ws.Range("M17").Value = 87.804878048780
$ws.Range("N17").Value = -23.762376237623

# ---------------------------------------------------------------------
# Row 18 - Burglary.
# ---------------------------------------------------------------------
Set-Count "C18" 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 3
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 48
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = 41.176470588235
$ws.Range("L18").Value = 118.181818181818
$ws.Range("M18").Value = -49.473684210526
$ws.Range("N18").Value = -83.892617449664

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny.
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 700
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 19
$ws.Range("H19").Value = 31.578947368421
$ws.Range("I19").Value = 249
$ws.Range("J19").Value = 233
$ws.Range("K19").Value = 6.866952789699
$ws.Range("L19").Value = 64.900662251655
$ws.Range("M19").Value = 99.2
$ws.Range("N19").Value = 19.711538461538

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -56.25
$ws.Range("I20").Value = 71
$ws.Range("J20").Value = 106
$ws.Range("K20").Value = -33.018867924528
$ws.Range("L20").Value = 44.897959183673
$ws.Range("M20").Value = 115.151515151515
$ws.Range("N20").Value = -88.854003139717

# ---------------------------------------------------------------------
# Row 21 - TOTAL.
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 100
$ws.Range("F21").Value = 41
$ws.Range("G21").Value = 47
$ws.Range("H21").Value = -12.765957446808
$ws.Range("I21").Value = 465
$ws.Range("J21").Value = 437
$ws.Range("K21").Value = 6.407322654462
$ws.Range("L21").Value = 76.806083650190
$ws.Range("M21").Value = 45.3125
$ws.Range("N21").Value = -64.422341239479

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny.
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -27.777777777777
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 49
$ws.Range("H24").Value = 2.040816326530
$ws.Range("I24").Value = 429
$ws.Range("J24").Value = 419
$ws.Range("K24").Value = 2.386634844868
$ws.Range("L24").Value = 86.521739130434
$ws.Range("M24").Value = -8.137044967880

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault.
# ---------------------------------------------------------------------
Set-Count "C25" 2
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -60
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = -54.545454545454
$ws.Range("I25").Value = 168
$ws.Range("J25").Value = 161
$ws.Range("K25").Value = 4.347826086956
$ws.Range("L25").Value = 38.842975206611
$ws.Range("M25").Value = -12.041884816753

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*.
# ---------------------------------------------------------------------
Set-TextZero "D26"
Set-TextStar "E26"

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes.
# ---------------------------------------------------------------------
Set-Count "D27" 1
Set-Pct "E27" -100 $fmtPct1
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = 60

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes.
# ---------------------------------------------------------------------
Set-TextZero "D30"
Set-TextStar "E30"
